# Auto-generated edit script: refresh crypto price/volume data and two pairs of
# rank-swapped rows, per commit "Updated cryptos list on Mon Dec 25 20:14:00 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.592.30'
$ws.Range("D3").Value = '2.271.03'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '121.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '265.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.644'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.30%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0944'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.911'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = '2.613.85'
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '2.265.71'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '43.545.79'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0915'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.75%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.16%  '
$ws.Range("E37").Value = '  +5.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.60'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.90%  '
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.25'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.87'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.237'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '73.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +40.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.34%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.100'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '101.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.64%  '
